$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: "U19" number format changes from date-only to date-time (style 3 -> style 2); value (2025-04-04) is unchanged
$ws.Range("U19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append new row 20 of bunker price data (dimension grows from A1:AV19 to A1:AV20)
$ws.Range("D20").Value2 = 552
$ws.Range("F20").Value2 = 560
$ws.Range("G20").Value2 = 613
$ws.Range("I20").Value2 = 543
$ws.Range("J20").Value2 = 546
$ws.Range("K20").Value2 = 540
$ws.Range("L20").Value2 = 580
$ws.Range("N20").Value2 = 515
$ws.Range("O20").Value2 = 501
$ws.Range("P20").Value2 = 654
$ws.Range("R20").Value2 = 546
$ws.Range("S20").Value2 = 523
$ws.Range("T20").Value2 = 567
$ws.Range("U20").Value2 = 45747
$ws.Range("X20").Value2 = 878
$ws.Range("Y20").Value2 = 537
$ws.Range("Z20").Value2 = 580
$ws.Range("AB20").Value2 = 590
$ws.Range("AC20").Value2 = 527
$ws.Range("AG20").Value2 = 674
$ws.Range("AH20").Value2 = 640
$ws.Range("AI20").Value2 = 630
$ws.Range("AJ20").Value2 = 501
$ws.Range("AK20").Value2 = 532
$ws.Range("AL20").Value2 = 609.5
$ws.Range("AM20").Value2 = 655
$ws.Range("AN20").Value2 = 624
$ws.Range("AP20").Value2 = 600
$ws.Range("AQ20").Value2 = 532
$ws.Range("AR20").Value2 = 530
$ws.Range("AS20").Value2 = 582
$ws.Range("AT20").Value2 = 656
$ws.Range("AV20").Value2 = 660

# "U20" (Date column) keeps the date-only format that "U19" used to have
$ws.Range("U20").NumberFormat = "YYYY-MM-DD"
